$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# NOTE: all row numbers below are the ORIGINAL (pre-edit) 1-indexed row
# positions in the table. We process the table from the bottom upward so
# that row numbers for not-yet-touched rows never shift underneath us.

# --- Last three summary rows: collapse the tab-separated run list down to
# a single value (rows 46, 45, 44 in that order, bottom-most first) ---
$t.Rows.Item(46).Cells.Item(1).Range.Text = "1432"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "1.62"
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.89"

# --- Row 12 (was 0.78518) becomes 1.62430 ---
$t.Rows.Item(12).Cells.Item(1).Range.Text = "1.62430"

# --- Remove rows 9, 10, 11 (0.00014 / 0.00016 / 0.00019), highest index
# first so the remaining deletions still target the right rows ---
$t.Rows.Item(11).Delete()
$t.Rows.Item(10).Delete()
$t.Rows.Item(9).Delete()

# --- Rows 6, 7, 8 each become 0.00558 ---
$t.Rows.Item(8).Cells.Item(1).Range.Text = "0.00558"
$t.Rows.Item(7).Cells.Item(1).Range.Text = "0.00558"
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00558"

# --- Insert three new rows right after row 4 (the one that becomes 7906),
# populated with 0.00002 / 0.00558 / 0.00107 in that order ---
$newRow1 = $t.Rows.Add($t.Rows.Item(5))
$newRow1.Cells.Item(1).Range.Text = "0.00002"
$newRow2 = $t.Rows.Add($t.Rows.Item(6))
$newRow2.Cells.Item(1).Range.Text = "0.00558"
$newRow3 = $t.Rows.Add($t.Rows.Item(7))
$newRow3.Cells.Item(1).Range.Text = "0.00107"

# --- Row 4 (was 5791) becomes 7906 ---
$t.Rows.Item(4).Cells.Item(1).Range.Text = "7906"

# --- Rows 1, 2, 3 (99.89 / 1.62 / 1432) all become 0M ---
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"

Write-Output ("Final row count=" + $t.Rows.Count)
